$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values to restore for each target row, pulled from the corresponding
# source row in the original (pre-edit) layout.
# Mapping: target row -> source row (both in the 2..8 range)
$rowData = @{
    2 = @{ D = 44208; L = 'Especial';     M = 70;  N = 24000; O = 24000; P = 24000; Q = '$/caja 15 kilos granel';    S = 1600; T = 15 }
    3 = @{ D = 44418; L = 'Especial';     M = 100; N = 8000;  O = 8000;  P = 8000;  Q = '$/caja 15 kilos granel';    S = 533;  T = 15 }
    4 = @{ D = 44392; L = 'Especial';     M = 500; N = 7000;  O = 7000;  P = 7000;  Q = '$/bandeja 8 kilos';         S = 875;  T = 8 }
    5 = @{ D = 44217; L = 'Primera';      M = 55;  N = 18000; O = 18000; P = 18000; Q = '$/caja 18 kilos granel';    S = 1000; T = 18 }
    6 = @{ D = 44264; L = 'Calibre 100';  M = 50;  N = 20000; O = 20000; P = 20000; Q = '$/caja 18 kilos embalada';  S = 1111; T = 18 }
    7 = @{ D = 44427; L = 'Primera';      M = 55;  N = 7000;  O = 7000;  P = 7000;  Q = '$/caja 15 kilos granel';    S = 467;  T = 15 }
    8 = @{ D = 44411; L = 'Primera';      M = 210; N = 8000;  O = 8000;  P = 8000;  Q = '$/bandeja 8 kilos';         S = 1000; T = 8 }
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("S$r").Value = $vals.S
    $ws.Range("T$r").Value = $vals.T
}
